$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 289
$ws.Range("A289").Value = 1
$ws.Range("B289").Value = 1
$ws.Range("C289").Value = "bozinovic2016"
$ws.Range("D289").Value = "Figure 2"
$ws.Range("E289").Value = 0
$ws.Range("F289").Value = 1
$ws.Range("H289").Value = 20
$ws.Range("I289").Value = 10
$ws.Range("J289").Value = 24
$ws.Range("K289").Value = "CTmin"
$ws.Range("L289").Value = "C"
$ws.Range("M289").Value = 7.9081911772020703
$ws.Range("N289").Value = 8.4562590142737406
$ws.Range("O289").Value = 0.28845675635350965
$ws.Range("P289").Value = 0.28845675635350965
$ws.Range("Q289").Value = 200
$ws.Range("R289").Value = 200
$ws.Range("S289").Value = 1
$ws.Range("T289").Value = 1
$ws.Range("U289").Value = "Drosophila "
$ws.Range("V289").Value = "melanogaster"
$ws.Range("W289").Value = 1
$ws.Range("X289").Value = 2
$ws.Range("Y289").Value = 1

# Row 290
$ws.Range("A290").Value = 2
$ws.Range("B290").Value = 1
$ws.Range("C290").Value = "bozinovic2016"
$ws.Range("D290").Value = "Figure 2"
$ws.Range("E290").Value = 0
$ws.Range("F290").Value = 1
$ws.Range("H290").Value = 20
$ws.Range("I290").Value = 10
$ws.Range("J290").Value = 24
$ws.Range("K290").Value = "CTmax"
$ws.Range("L290").Value = "C"
$ws.Range("M290").Value = 37.837965140000001
$ws.Range("M290").Font.Color = 0
$ws.Range("N290").Value = 38.114178049929301
$ws.Range("O290").Value = 0.21818181818180449
$ws.Range("P290").Value = 0.21818181818179738
$ws.Range("Q290").Value = 200
$ws.Range("R290").Value = 200
$ws.Range("S290").Value = 1
$ws.Range("T290").Value = 1
$ws.Range("U290").Value = "Drosophila "
$ws.Range("V290").Value = "melanogaster"
$ws.Range("W290").Value = 1
$ws.Range("X290").Value = 2
$ws.Range("Y290").Value = 1

# Row 291
$ws.Range("A291").Value = 1
$ws.Range("B291").Value = 2
$ws.Range("C291").Value = "bozinovic2016"
$ws.Range("D291").Value = "Figure 3"
$ws.Range("E291").Value = 0
$ws.Range("F291").Value = 1
$ws.Range("H291").Value = 20
$ws.Range("I291").Value = 10
$ws.Range("J291").Value = 24
$ws.Range("K291").Value = "Scope of thermal tolerance (CTmax - CTmin)"
$ws.Range("L291").Value = "C"
$ws.Range("M291").Value = 29.868428605139101
$ws.Range("N291").Value = 29.594316217784701
$ws.Range("O291").Value = 0.24911190902889757
$ws.Range("P291").Value = 0.24911190902889757
$ws.Range("Q291").Value = 200
$ws.Range("R291").Value = 200
$ws.Range("S291").Value = 1
$ws.Range("T291").Value = 1
$ws.Range("U291").Value = "Drosophila "
$ws.Range("V291").Value = "melanogaster"
$ws.Range("W291").Value = 1
$ws.Range("X291").Value = 2
$ws.Range("Y291").Value = 1

# Row 292
$ws.Range("A292").Value = 1
$ws.Range("B292").Value = 1
$ws.Range("C292").Value = "delnat2019"
$ws.Range("D292").Value = "Figure 2a"
$ws.Range("E292").Value = 0
$ws.Range("F292").Value = 1
$ws.Range("H292").Value = 20
$ws.Range("I292").Value = 7
$ws.Range("J292").Value = 24
$ws.Range("K292").Value = "CTmax"
$ws.Range("L292").Value = "C"
$ws.Range("M292").Value = 40.6
$ws.Range("N292").Value = 40.864150943396197
$ws.Range("O292").Value = 0.30377358490564887
$ws.Range("P292").Value = 0.30377358490570217
$ws.Range("Q292").Value = 75
$ws.Range("R292").Value = 78
$ws.Range("S292").Value = 1
$ws.Range("T292").Value = 0
$ws.Range("U292").Value = "Culex"
$ws.Range("V292").Value = "pipiens"
$ws.Range("W292").Value = 1
$ws.Range("X292").Value = 0
$ws.Range("Y292").Value = 1
$ws.Range("Z292").Value = "insecticide exposure"
$ws.Range("AA292").Value = "absence"

# Row 293
$ws.Range("A293").Value = 1
$ws.Range("B293").Value = 1
$ws.Range("C293").Value = "delnat2019"
$ws.Range("D293").Value = "Figure 2a"
$ws.Range("E293").Value = 0
$ws.Range("F293").Value = 1
$ws.Range("H293").Value = 20
$ws.Range("I293").Value = 7
$ws.Range("J293").Value = 24
$ws.Range("K293").Value = "CTmax"
$ws.Range("L293").Value = "C"
$ws.Range("M293").Value = 36.796226415094303
$ws.Range("N293").Value = 35.9509433962264
$ws.Range("O293").Value = 0.30377358490569861
$ws.Range("P293").Value = 0.31698113207544765
$ws.Range("Q293").Value = 78
$ws.Range("R293").Value = 76
$ws.Range("S293").Value = 1
$ws.Range("T293").Value = 0
$ws.Range("U293").Value = "Culex"
$ws.Range("V293").Value = "pipiens"
$ws.Range("W293").Value = 1
$ws.Range("X293").Value = 0
$ws.Range("Y293").Value = 1
$ws.Range("Z293").Value = "insecticide exposure"
$ws.Range("AA293").Value = "presence"

# Row 294
$ws.Range("A294").Value = 1
$ws.Range("B294").Value = 1
$ws.Range("C294").Value = "delnat2019"
$ws.Range("D294").Value = "Figure 2a"
$ws.Range("E294").Value = 0
$ws.Range("F294").Value = 1
$ws.Range("H294").Value = 20
$ws.Range("I294").Value = 14
$ws.Range("J294").Value = 24
$ws.Range("K294").Value = "CTmax"
$ws.Range("L294").Value = "C"
$ws.Range("M294").Value = 40.6
$ws.Range("N294").Value = 41.022641509433903
$ws.Range("O294").Value = 0.30377358490564887
$ws.Range("P294").Value = 0.31698113207550094
$ws.Range("Q294").Value = 75
$ws.Range("R294").Value = 72
$ws.Range("S294").Value = 1
$ws.Range("T294").Value = 0
$ws.Range("U294").Value = "Culex"
$ws.Range("V294").Value = "pipiens"
$ws.Range("W294").Value = 1
$ws.Range("X294").Value = 0
$ws.Range("Y294").Value = 1
$ws.Range("Z294").Value = "insecticide exposure"
$ws.Range("AA294").Value = "absence"

# Row 295
$ws.Range("A295").Value = 1
$ws.Range("B295").Value = 1
$ws.Range("C295").Value = "delnat2019"
$ws.Range("D295").Value = "Figure 2a"
$ws.Range("E295").Value = 0
$ws.Range("F295").Value = 1
$ws.Range("H295").Value = 20
$ws.Range("I295").Value = 14
$ws.Range("J295").Value = 24
$ws.Range("K295").Value = "CTmax"
$ws.Range("L295").Value = "C"
$ws.Range("M295").Value = 36.796226415094303
$ws.Range("N295").Value = 38.090566037735798
$ws.Range("O295").Value = 0.30377358490569861
$ws.Range("P295").Value = 0.34339622641509848
$ws.Range("Q295").Value = 78
$ws.Range("R295").Value = 68
$ws.Range("S295").Value = 1
$ws.Range("T295").Value = 0
$ws.Range("U295").Value = "Culex"
$ws.Range("V295").Value = "pipiens"
$ws.Range("W295").Value = 1
$ws.Range("X295").Value = 0
$ws.Range("Y295").Value = 1
$ws.Range("Z295").Value = "insecticide exposure"
$ws.Range("AA295").Value = "presence"

# Row 296
$ws.Range("A296").Value = 1
$ws.Range("B296").Value = 2
$ws.Range("C296").Value = "delnat2019"
$ws.Range("D296").Value = "Figure 2b"
$ws.Range("E296").Value = 0
$ws.Range("F296").Value = 1
$ws.Range("H296").Value = 20
$ws.Range("I296").Value = 7
$ws.Range("J296").Value = 24
$ws.Range("K296").Value = "CTmax"
$ws.Range("L296").Value = "C"
$ws.Range("M296").Value = 39.411320754716897
$ws.Range("N296").Value = 40.098113207547101
$ws.Range("O296").Value = 0.27735849056604778
$ws.Range("P296").Value = 0.33018867924524997
$ws.Range("Q296").Value = 106
$ws.Range("R296").Value = 76
$ws.Range("S296").Value = 1
$ws.Range("T296").Value = 0
$ws.Range("U296").Value = "Culex"
$ws.Range("V296").Value = "pipiens"
$ws.Range("W296").Value = 1
$ws.Range("X296").Value = 2
$ws.Range("Y296").Value = 1
$ws.Range("Z296").Value = "insecticide exposure"
$ws.Range("AA296").Value = "absence"
$ws.Range("AB296").Value = "sex"
$ws.Range("AC296").Value = "male"

# Row 297
$ws.Range("A297").Value = 1
$ws.Range("B297").Value = 2
$ws.Range("C297").Value = "delnat2019"
$ws.Range("D297").Value = "Figure 2b"
$ws.Range("E297").Value = 0
$ws.Range("F297").Value = 1
$ws.Range("H297").Value = 20
$ws.Range("I297").Value = 7
$ws.Range("J297").Value = 24
$ws.Range("K297").Value = "CTmax"
$ws.Range("L297").Value = "C"
$ws.Range("M297").Value = 36.2415094339622
$ws.Range("N297").Value = 37.7735849056603
$ws.Range("O297").Value = 0.40943396226414919
$ws.Range("P297").Value = 0.36981132075474932
$ws.Range("Q297").Value = 50
$ws.Range("R297").Value = 60
$ws.Range("S297").Value = 1
$ws.Range("T297").Value = 0
$ws.Range("U297").Value = "Culex"
$ws.Range("V297").Value = "pipiens"
$ws.Range("W297").Value = 1
$ws.Range("X297").Value = 2
$ws.Range("Y297").Value = 1
$ws.Range("Z297").Value = "insecticide exposure"
$ws.Range("AA297").Value = "presence"
$ws.Range("AB297").Value = "sex"
$ws.Range("AC297").Value = "male"

# Row 298
$ws.Range("A298").Value = 1
$ws.Range("B298").Value = 2
$ws.Range("C298").Value = "delnat2019"
$ws.Range("D298").Value = "Figure 2b"
$ws.Range("E298").Value = 0
$ws.Range("F298").Value = 1
$ws.Range("H298").Value = 20
$ws.Range("I298").Value = 14
$ws.Range("J298").Value = 24
$ws.Range("K298").Value = "CTmax"
$ws.Range("L298").Value = "C"
$ws.Range("M298").Value = 39.411320754716897
$ws.Range("N298").Value = 40.071698113207503
$ws.Range("O298").Value = 0.27735849056604778
$ws.Range("P298").Value = 0.30377358490564887
$ws.Range("Q298").Value = 106
$ws.Range("R298").Value = 89
$ws.Range("S298").Value = 1
$ws.Range("T298").Value = 0
$ws.Range("U298").Value = "Culex"
$ws.Range("V298").Value = "pipiens"
$ws.Range("W298").Value = 1
$ws.Range("X298").Value = 2
$ws.Range("Y298").Value = 1
$ws.Range("Z298").Value = "insecticide exposure"
$ws.Range("AA298").Value = "absence"
$ws.Range("AB298").Value = "sex"
$ws.Range("AC298").Value = "male"

# Row 299
$ws.Range("A299").Value = 1
$ws.Range("B299").Value = 2
$ws.Range("C299").Value = "delnat2019"
$ws.Range("D299").Value = "Figure 2b"
$ws.Range("E299").Value = 0
$ws.Range("F299").Value = 1
$ws.Range("H299").Value = 20
$ws.Range("I299").Value = 14
$ws.Range("J299").Value = 24
$ws.Range("K299").Value = "CTmax"
$ws.Range("L299").Value = "C"
$ws.Range("M299").Value = 36.2415094339622
$ws.Range("N299").Value = 37.562264150943399
$ws.Range("O299").Value = 0.40943396226414919
$ws.Range("P299").Value = 0.44905660377354906
$ws.Range("Q299").Value = 50
$ws.Range("R299").Value = 44
$ws.Range("S299").Value = 1
$ws.Range("T299").Value = 0
$ws.Range("U299").Value = "Culex"
$ws.Range("V299").Value = "pipiens"
$ws.Range("W299").Value = 1
$ws.Range("X299").Value = 2
$ws.Range("Y299").Value = 1
$ws.Range("Z299").Value = "insecticide exposure"
$ws.Range("AA299").Value = "presence"
$ws.Range("AB299").Value = "sex"
$ws.Range("AC299").Value = "male"

# Row 300
$ws.Range("A300").Value = 1
$ws.Range("B300").Value = 3
$ws.Range("C300").Value = "delnat2019"
$ws.Range("D300").Value = "Figure 2c"
$ws.Range("E300").Value = 0
$ws.Range("F300").Value = 1
$ws.Range("H300").Value = 20
$ws.Range("I300").Value = 14
$ws.Range("J300").Value = 24
$ws.Range("K300").Value = "CTmax"
$ws.Range("L300").Value = "C"
$ws.Range("M300").Value = 40.018867924528301
$ws.Range("N300").Value = 40.784905660377298
$ws.Range("O300").Value = 0.36981132075474932
$ws.Range("P300").Value = 0.39622641509435041
$ws.Range("Q300").Value = 84
$ws.Range("R300").Value = 63
$ws.Range("S300").Value = 1
$ws.Range("T300").Value = 0
$ws.Range("U300").Value = "Culex"
$ws.Range("V300").Value = "pipiens"
$ws.Range("W300").Value = 1
$ws.Range("X300").Value = 2
$ws.Range("Y300").Value = 1
$ws.Range("Z300").Value = "insecticide exposure"
$ws.Range("AA300").Value = "presence"
$ws.Range("AB300").Value = "sex"
$ws.Range("AC300").Value = "female"

# Row 301
$ws.Range("A301").Value = 1
$ws.Range("B301").Value = 3
$ws.Range("C301").Value = "delnat2019"
$ws.Range("D301").Value = "Figure 2c"
$ws.Range("E301").Value = 0
$ws.Range("F301").Value = 1
$ws.Range("H301").Value = 20
$ws.Range("I301").Value = 14
$ws.Range("J301").Value = 24
$ws.Range("K301").Value = "CTmax"
$ws.Range("L301").Value = "C"
$ws.Range("M301").Value = 39.543396226415098
$ws.Range("N301").Value = 38.539622641509403
$ws.Range("O301").Value = 0.38301886792454809
$ws.Range("P301").Value = 0.39622641509435041
$ws.Range("Q301").Value = 68
$ws.Range("R301").Value = 60
$ws.Range("S301").Value = 1
$ws.Range("T301").Value = 0
$ws.Range("U301").Value = "Culex"
$ws.Range("V301").Value = "pipiens"
$ws.Range("W301").Value = 1
$ws.Range("X301").Value = 2
$ws.Range("Y301").Value = 1
$ws.Range("Z301").Value = "insecticide exposure"
$ws.Range("AA301").Value = "presence"
$ws.Range("AB301").Value = "sex"
$ws.Range("AC301").Value = "female"

# Row 302
$ws.Range("A302").Value = 1
$ws.Range("B302").Value = 3
$ws.Range("C302").Value = "delnat2019"
$ws.Range("D302").Value = "Figure 2c"
$ws.Range("E302").Value = 0
$ws.Range("F302").Value = 1
$ws.Range("H302").Value = 20
$ws.Range("I302").Value = 14
$ws.Range("J302").Value = 24
$ws.Range("K302").Value = "CTmax"
$ws.Range("L302").Value = "C"
$ws.Range("M302").Value = 40.018867924528301
$ws.Range("N302").Value = 41.075471698113198
$ws.Range("O302").Value = 0.36981132075474932
$ws.Range("P302").Value = 0.42264150943395151
$ws.Range("Q302").Value = 84
$ws.Range("R302").Value = 57
$ws.Range("S302").Value = 1
$ws.Range("T302").Value = 0
$ws.Range("U302").Value = "Culex"
$ws.Range("V302").Value = "pipiens"
$ws.Range("W302").Value = 1
$ws.Range("X302").Value = 2
$ws.Range("Y302").Value = 1
$ws.Range("Z302").Value = "insecticide exposure"
$ws.Range("AA302").Value = "presence"
$ws.Range("AB302").Value = "sex"
$ws.Range("AC302").Value = "female"

# Row 303
$ws.Range("A303").Value = 1
$ws.Range("B303").Value = 3
$ws.Range("C303").Value = "delnat2019"
$ws.Range("D303").Value = "Figure 2c"
$ws.Range("E303").Value = 0
$ws.Range("F303").Value = 1
$ws.Range("H303").Value = 20
$ws.Range("I303").Value = 14
$ws.Range("J303").Value = 24
$ws.Range("K303").Value = "CTmax"
$ws.Range("L303").Value = "C"
$ws.Range("M303").Value = 39.543396226415098
$ws.Range("N303").Value = 38.803773584905599
$ws.Range("O303").Value = 0.38301886792454809
$ws.Range("P303").Value = 0.38301886792455164
$ws.Range("Q303").Value = 68
$ws.Range("R303").Value = 65
$ws.Range("S303").Value = 1
$ws.Range("T303").Value = 0
$ws.Range("U303").Value = "Culex"
$ws.Range("V303").Value = "pipiens"
$ws.Range("W303").Value = 1
$ws.Range("X303").Value = 2
$ws.Range("Y303").Value = 1
$ws.Range("Z303").Value = "insecticide exposure"
$ws.Range("AA303").Value = "presence"
$ws.Range("AB303").Value = "sex"
$ws.Range("AC303").Value = "female"
